# massive MDY site template update
# Insert three new columns (Month, Day, Year) before the existing
# "Date Sampled" column on the "Data Entry" sheet, populate them from the
# existing Date Sampled values, and restore the prior formatting /
# selection state as closely as possible.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Entry")

# Insert 3 new columns at E:G (old column E "Date Sampled" and everything
# to its right shifts right by 3 columns). Excel inherits formatting from
# the column immediately to the left (D) for the newly inserted columns,
# which is exactly what the target workbook shows (header style matches
# D1, data-row style matches D2:D5).
$ws.Range("E1:G1").EntireColumn.Insert()

# Give the three new columns the same column width as column D, mirroring
# the original author's "insert columns" workflow.
$ws.Range("E1:G1").EntireColumn.ColumnWidth = $ws.Range("D1").EntireColumn.ColumnWidth

# New header labels
$ws.Range("E1").Value = "Month"
$ws.Range("F1").Value = "Day"
$ws.Range("G1").Value = "Year"

# New Month/Day/Year values, derived from the existing Date Sampled column
# (now shifted to column H).
$lastRow = 5
for ($r = 2; $r -le $lastRow; $r++) {
    $dateCell = $ws.Cells.Item($r, 8)
    $d = $dateCell.Value2
    $dt = [DateTime]::FromOADate($d)
    $ws.Cells.Item($r, 5).Value = $dt.Month
    $ws.Cells.Item($r, 6).Value = $dt.Day
    $ws.Cells.Item($r, 7).Value = $dt.Year
}

# Restore the active selection shown in the saved workbook
$ws.Range("E6").Select() | Out-Null
